# HP: Operation error message propagation.
# Update the "Results" sheet's result messages to include the "RESULT: " / "ERROR: "
# prefixes used for propagating operation errors, and surface the operation
# error text for the 4th submission in a new H14 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("G11").Value = "RESULT: SUBMIT(1) TESTCASE(1) MSG:(TESTCASE#1:EQUAL`n)"
$ws.Range("G12").Value = "RESULT: SUBMIT(2) TESTCASE(1) MSG:(TESTCASE#1:NOT_EQUAL`n)"
$ws.Range("G13").Value = "RESULT: SUBMIT(3) TESTCASE(1) MSG:(TESTCASE#1:NOT_EQUAL`n)"
$ws.Range("G14").Value = "RESULT: SUBMIT(4) TESTCASE(1) MSG:(TESTCASE#1:NOT_EQUAL`n)"
$ws.Range("H14").Value = "ERROR: SUBMIT(4) TESTCASE(1) MSG:(CLASS:siima.app.XSLTransformer ERROR:Syntax error in '/CATALOG/Plant[(Price>'6.60')] and [(Light='Sun')]'.)"

# Setting the wrapped-text G-column cells above triggers Excel's automatic
# row-height autofit; re-running AutoFit on these (already-default-height)
# rows restores them to the sheet's default row height so no stray
# per-row height override is left behind.
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
